# Adds Statistical analyses ("ps") species rows to the SA code dataset,
# plus Surface Area values for the trials that have raw scale measurements.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49 is the last fully-populated template row; reuse its number
# formatting/styling for every new row we are about to populate.
$ws.Range("A49:K49").Copy()

$rowsData = @(
    @{Row=50; A=4.0; C=7.0; D=3.0; E="wild"; F=0.761492; G=0.766028; H=0.756092; K=190.016739},
    @{Row=51; A=4.0; C=11.0; D=2.0; E="wild"; F=0.916156; G=0.916612; H=0.901219; K=103.633705},
    @{Row=52; A=4.0; C=12.0; D=2.0; E="wild"; F=1.54855; G=1.54116; H=1.53285; K=177.580276},
    @{Row=53; A=4.0; C=13.0; D=1.0; E="wild"; F=1.56975; G=1.56975; H=1.56975; K=252.151855},
    @{Row=54; A=4.0; C=14.0; D=1.0; E="wild"; F=6.5366; G=6.51526; H=6.49434; K=280.506165},
    @{Row=55; A=4.0; C=7.0; D=3.0; E="frag"; F=0.478508; G=0.482802; H=0.48883; K=135.365265},
    @{Row=56; A=4.0; C=11.0; D=3.0; E="frag"; F=1.82289; G=1.84345; H=1.87992; K=34.749279},
    @{Row=57; A=4.0; C=12.0; D=3.0; E="frag"; F=1.25343; G=1.26447; H=1.25494; K=87.473793},
    @{Row=58; A=4.0; C=13.0; D=3.0; E="frag"; F=1.33362; G=1.33067; H=1.34269; K=140.08075},
    @{Row=59; A=4.0; C=14.0; D=3.0; E="frag"; F=1.07647; G=1.07574; H=1.08931; K=191.489624},
    @{Row=60; A=4.0; C=7.0; D=2.0; E="frag"; F=2.93244; G=2.89922; H=2.92458},
    @{Row=61; A=4.0; C=11.0; D=2.0; E="frag"},
    @{Row=62; A=4.0; C=12.0; D=2.0; E="frag"},
    @{Row=63; A=4.0; C=13.0; D=2.0; E="frag"},
    @{Row=64; A=4.0; C=14.0; D=2.0; E="frag"},
    @{Row=65; A=4.0; C=7.0; D=1.0; E="frag"},
    @{Row=66; A=4.0; C=11.0; D=1.0; E="frag"},
    @{Row=67; A=4.0; C=12.0; D=1.0; E="frag"},
    @{Row=68; A=4.0; C=13.0; D=1.0; E="frag"},
    @{Row=69; A=4.0; C=14.0; D=1.0; E="frag"},
    @{Row=70; A=0.0; C=7.0; D=3.0; E="wild"},
    @{Row=71; A=0.0; C=11.0; D=2.0; E="wild"},
    @{Row=72; A=0.0; C=12.0; D=2.0; E="wild"},
    @{Row=73; A=0.0; C=13.0; D=1.0; E="wild"},
    @{Row=74; A=0.0; C=14.0; D=1.0; E="wild"},
    @{Row=75; A=0.0; C=7.0; D=3.0; E="frag"},
    @{Row=76; A=0.0; C=11.0; D=3.0; E="frag"},
    @{Row=77; A=0.0; C=12.0; D=3.0; E="frag"},
    @{Row=78; A=0.0; C=13.0; D=3.0; E="frag"},
    @{Row=79; A=0.0; C=14.0; D=3.0; E="frag"},
    @{Row=80; A=0.0; C=7.0; D=2.0; E="frag"},
    @{Row=81; A=0.0; C=11.0; D=2.0; E="frag"},
    @{Row=82; A=0.0; C=12.0; D=2.0; E="frag"},
    @{Row=83; A=0.0; C=13.0; D=2.0; E="frag"},
    @{Row=84; A=0.0; C=14.0; D=2.0; E="frag"},
    @{Row=85; A=0.0; C=7.0; D=1.0; E="frag"},
    @{Row=86; A=0.0; C=11.0; D=1.0; E="frag"},
    @{Row=87; A=0.0; C=12.0; D=1.0; E="frag"},
    @{Row=88; A=0.0; C=13.0; D=1.0; E="frag"},
    @{Row=89; A=0.0; C=14.0; D=1.0; E="frag"}

)

foreach ($rd in $rowsData) {
    $row = $rd.Row

    if ($rd.ContainsKey("K")) {
        # full data row: A-E identifiers, F-H raw scale reads, K surface area
        $ws.Range("A$row`:K$row").PasteSpecial(-4122)
    } elseif ($rd.ContainsKey("F")) {
        # data row with raw scale reads but no computed surface area yet
        $ws.Range("A$row`:J$row").PasteSpecial(-4122)
    } else {
        # identifier-only row (no measurements recorded)
        $ws.Range("A$row`:E$row").PasteSpecial(-4122)
    }

    $ws.Cells.Item($row, 1).Value = $rd.A          # A: timePoint
    $ws.Cells.Item($row, 2).Value = "ps"            # B: species
    $ws.Cells.Item($row, 3).Value = $rd.C          # C: genotype
    $ws.Cells.Item($row, 4).Value = $rd.D          # D: trial
    $ws.Cells.Item($row, 5).Value = $rd.E          # E: ID

    if ($rd.ContainsKey("F")) {
        $ws.Cells.Item($row, 6).Value = $rd.F      # F: Scale 1
        $ws.Cells.Item($row, 7).Value = $rd.G      # G: Scale 2
        $ws.Cells.Item($row, 8).Value = $rd.H      # H: Scale 3
    }

    if ($rd.ContainsKey("K")) {
        $ws.Cells.Item($row, 11).Value = $rd.K     # K: Surface Area (cm2)
    }
}

# Row 90 gets only the styled (empty) timePoint cell in column A.
$ws.Range("A49").Copy()
$ws.Range("A90").PasteSpecial(-4122)

# Extend the Scale Factor formula down to the newly added rows (J2:J83 -> J2:J90).
# Copy J83's number formatting first, then write the shared formula across
# J84:J90 in one assignment so Excel groups it as a single shared formula.
$ws.Range("J83").Copy()
$ws.Range("J84:J90").PasteSpecial(-4122)
$ws.Range("J84:J90").Formula = "=3/I84"

Write-Host "Applied SA statistical-analysis rows 50-90"
